# Capitalization of column headings in tables
# (plus: title gets " RRF" suffix, one header label simplified, font sizes
#  bumped from 11pt to 12pt, selection + row-1 height follow the edit)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Text content changes -------------------------------------------------
# Title row (merged B1:I1)
$ws.Range("B1").Value = "LLM modell: gemini-2.0-flash;  RRF"

# Column headings (row 2) - capitalised, one re-worded
$ws.Range("B2").Value = "Kérdések száma"
$ws.Range("C2").Value = " Embedding  generálásai idő átlaga"
$ws.Range("D2").Value = "Sparse embedding generálási idő átlaga"
$ws.Range("E2").Value = "Kontextus összeállitási idő átlaga"
$ws.Range("G2").Value = "Teljes feldoldozási idő átlaga"
$ws.Range("H2").Value = "Szemantikus hasonlóság mérékének  (BERTScore F1) átlaga (0-1) között"
$ws.Range("I2").Value = "Top_k darab számának átlag"

# ---- Font size bump: 11pt -> 12pt ------------------------------------------
$ws.Range("B1").Font.Size = 12
$ws.Range("C1:I1").Font.Size = 12
$ws.Range("B2:I2").Font.Size = 12
$ws.Range("B3:I5").Font.Size = 12

# ---- Row height & selection -------------------------------------------------
$ws.Rows.Item(1).RowHeight = 15.6
[void]$ws.Range("E11").Select()
